$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 58
$ws.Range("H58").Value = 2100.8
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 2100.8
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 6302.400000000001
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -6602.400000000001
# Row 80
$ws.Range("H80").Value = 941.5
$ws.Range("I80").Value = 575
$ws.Range("K80").Value = 1725
$ws.Range("M80").Value = -727
# Row 83
$ws.Range("H83").Value = 941.5
$ws.Range("I83").Value = 575
$ws.Range("K83").Value = 5175
$ws.Range("M83").Value = -183
# Row 92
$ws.Range("H92").Value = 177.375
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
# Row 98
$ws.Range("H98").Value = 8193.5
$ws.Range("I98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("M98").ClearContents()
# Row 122
$ws.Range("H122").Value = 8193.5
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
# Row 135
$ws.Range("H135").Value = 2218.7144
$ws.Range("I135").Value = 2218.7144
$ws.Range("K135").Value = 19968.4296
$ws.Range("M135").Value = -17433.4296
# Row 137
$ws.Range("H137").Value = 5275
$ws.Range("I137").Value = 3787.4
$ws.Range("K137").Value = 11362.2
$ws.Range("M137").Value = -8812.200000000001
# Row 138
$ws.Range("H138").Value = 2838.8
$ws.Range("I138").Value = 2462.6667
$ws.Range("J138").Value = 3000
$ws.Range("K138").Value = 7388.000100000001
$ws.Range("L138").Value = 9000
$ws.Range("M138").Value = -2248.000100000001
$ws.Range("N138").Value = -19280

$ws = $wb.Worksheets.Item("ARM")
# Row 95
$ws.Range("H95").Value = 17051.75
$ws.Range("J95").Value = 17051.75
$ws.Range("L95").Value = 17051.75
$ws.Range("N95").Value = -22543.75
# Row 96
$ws.Range("H96").Value = 33043.125
$ws.Range("J96").Value = 33043.125
$ws.Range("L96").Value = 33043.125
$ws.Range("N96").Value = -38535.125
# Row 110
$ws.Range("H110").Value = 249.5
$ws.Range("I110").Value = 199
$ws.Range("K110").Value = 199
$ws.Range("M110").Value = 1846
# Row 122
$ws.Range("H122").Value = 1488.9231
$ws.Range("I122").Value = 1488.9231
$ws.Range("K122").Value = 4466.7693
$ws.Range("M122").Value = -2016.7693

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 2686.6365
$ws.Range("I86").Value = 2962.875
$ws.Range("J86").Value = 1950
$ws.Range("K86").Value = 2962.875
$ws.Range("L86").Value = 1950
$ws.Range("M86").Value = -1839.875
$ws.Range("N86").Value = -4196
# Row 89
$ws.Range("H89").Value = 2686.6365
$ws.Range("I89").Value = 2962.875
$ws.Range("J89").Value = 1950
$ws.Range("K89").Value = 14814.375
$ws.Range("L89").Value = 9750
$ws.Range("M89").Value = -9198.375
$ws.Range("N89").Value = -20982
# Row 107
$ws.Range("H107").Value = 1099.5
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 1099.5
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 1099.5
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -4939.5
# Row 112
$ws.Range("H112").Value = 80000
$ws.Range("J112").Value = 80000
$ws.Range("L112").Value = 80000
$ws.Range("N112").Value = -82954

$ws = $wb.Worksheets.Item("CRP")
# Row 28
$ws.Range("H28").Value = 29549.834
$ws.Range("J28").Value = 29549.834
$ws.Range("L28").Value = 29549.834
$ws.Range("N28").Value = -30039.834
# Row 58
$ws.Range("H58").Value = 1778.6666
$ws.Range("I58").Value = 870.1111
$ws.Range("K58").Value = 870.1111
$ws.Range("M58").Value = -667.1111
# Row 107
$ws.Range("H107").Value = 1067.9286
$ws.Range("I107").Value = 1082.6666
$ws.Range("J107").Value = 1041.4
$ws.Range("K107").Value = 1082.6666
$ws.Range("L107").Value = 1041.4
$ws.Range("M107").Value = 837.3334
$ws.Range("N107").Value = -4881.4
# Row 136
$ws.Range("H136").Value = 1778.6666
$ws.Range("I136").Value = 870.1111
$ws.Range("K136").Value = 2610.3333
$ws.Range("M136").Value = -60.33329999999978

$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 180.85715
$ws.Range("J12").Value = 286.5
$ws.Range("L12").Value = 859.5
$ws.Range("N12").Value = -1205.5
# Row 108
$ws.Range("H108").Value = 1810.7142
$ws.Range("I108").Value = 1810.7142
$ws.Range("K108").Value = 5432.142599999999
$ws.Range("M108").Value = -2552.142599999999
# Row 119
$ws.Range("H119").Value = 9361.799999999999
$ws.Range("I119").Value = 7872.25
$ws.Range("K119").Value = 23616.75
$ws.Range("M119").Value = -18778.75
# Row 131
$ws.Range("H131").Value = 1610.909
$ws.Range("I131").Value = 1286.6666
$ws.Range("J131").Value = 2000
$ws.Range("K131").Value = 3859.9998
$ws.Range("L131").Value = 6000
$ws.Range("M131").Value = 1180.0002
$ws.Range("N131").Value = -16080

$ws = $wb.Worksheets.Item("GSM")
# Row 39
$ws.Range("H39").Value = 35000
$ws.Range("J39").Value = 35000
$ws.Range("L39").Value = 35000
$ws.Range("N39").Value = -36064
# Row 55
$ws.Range("H55").Value = 8000
$ws.Range("J55").Value = 8000
$ws.Range("L55").Value = 8000
$ws.Range("N55").Value = -8654
# Row 113
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 665
$ws.Range("I22").Value = 497.5
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 497.5
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -202.5
$ws.Range("N22").Value = -1590
# Row 27
$ws.Range("H27").Value = 665
$ws.Range("I27").Value = 497.5
$ws.Range("J27").Value = 1000
$ws.Range("K27").Value = 497.5
$ws.Range("L27").Value = 1000
$ws.Range("M27").Value = -390.5
$ws.Range("N27").Value = -1214
# Row 42
$ws.Range("H42").Value = 21000000
$ws.Range("I42").Value = 2000000
$ws.Range("K42").Value = 2000000
$ws.Range("M42").Value = -1999437
# Row 49
$ws.Range("H49").Value = 21000000
$ws.Range("I49").Value = 2000000
$ws.Range("K49").Value = 2000000
$ws.Range("M49").Value = -1999853
# Row 82
$ws.Range("H82").Value = 2004.5714
$ws.Range("J82").Value = 1922.5
$ws.Range("L82").Value = 1922.5
$ws.Range("N82").Value = -2644.5
# Row 85
$ws.Range("H85").Value = 2004.5714
$ws.Range("J85").Value = 1922.5
$ws.Range("L85").Value = 1922.5
$ws.Range("N85").Value = -4418.5

$ws = $wb.Worksheets.Item("WVR")
# Row 103
$ws.Range("H103").Value = 25000
$ws.Range("J103").Value = 25000
$ws.Range("L103").Value = 25000
$ws.Range("N103").Value = -27344
# Row 104
$ws.Range("H104").Value = 38900
$ws.Range("J104").Value = 38900
$ws.Range("L104").Value = 38900
$ws.Range("N104").Value = -45888
# Row 132
$ws.Range("H132").Value = 1840.8334
$ws.Range("I132").Value = 1839.25
$ws.Range("J132").Value = 1844
$ws.Range("K132").Value = 5517.75
$ws.Range("L132").Value = 5532
$ws.Range("M132").Value = -2987.75
$ws.Range("N132").Value = -10592
# Row 138
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
